$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value2 = [double]"7.720153894239088e-16"
$ws.Range("B2").Value2 = [double]"-0.9999999999999989"
$ws.Range("C2").Value2 = [double]"3.973817292551296e-08"
$ws.Range("D2").Value2 = [double]"-2.649211528367554e-08"

$ws.Range("A3").Value2 = [double]"1.58952691702053e-08"
$ws.Range("B3").Value2 = [double]"-3.973817292551311e-08"
$ws.Range("C3").Value2 = [double]"-0.9999999999999991"
$ws.Range("D3").Value2 = [double]"2.388943118032592e-15"

$ws.Range("A4").Value2 = [double]"1.324605764183779e-08"
$ws.Range("B4").Value2 = [double]"-2.649211528367546e-08"
$ws.Range("C4").Value2 = [double]"4.330927444391978e-15"
$ws.Range("D4").Value2 = [double]"0.9999999999999996"

$ws.Range("A5").Value2 = [double]"0.9999999999999998"
$ws.Range("B5").Value2 = [double]"1.754580430508887e-15"
$ws.Range("C5").Value2 = [double]"1.589526917020529e-08"
$ws.Range("D5").Value2 = [double]"-1.324605764183777e-08"
